# fix issues to get 'Cotação base' and 'Data de pagamento'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header columns (K1, L1), matching the style of the
# existing header cells (copy format from J1, the last existing header).
$ws.Range("K1").Value = "Próxima Cotação base"
$ws.Range("L1").Value = "Próxima data de pagamento"
$ws.Range("J1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)

# The old H2 ("Ultima data de pagamento" value, 75,00) actually belongs
# under the new K2 ("Próxima Cotação base") column - move it there.
$ws.Range("K2").Value = "75,00"

# Correct the G2/H2 values (Ultima Cotação base / Ultima data de pagamento)
# and fill in the new L2 (Próxima data de pagamento) value.
$ws.Range("G2").Value = "67,59"
$ws.Range("H2").Value = "15/05/2023"
$ws.Range("L2").Value = "15/06/2023"
